$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "datos actualizados" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Agosto de 2020 a las 19:45"

# --- Rows whose rank/position changed: country label + its stats move together ---
$ws.Cells.Item(68, 1).Value = "Etiopia"
$ws.Cells.Item(68, "B").Value = 25118
$ws.Cells.Item(68, "C").Value = 943
$ws.Cells.Item(68, "D").Value = 11034
$ws.Cells.Item(68, "E").Value = 13621
$ws.Cells.Item(68, "F").Value = 0
$ws.Cells.Item(68, "G").Value = 23
$ws.Cells.Item(68, "H").Value = 463

$ws.Cells.Item(69, 1).Value = "Costa Rica"
$ws.Cells.Item(69, "B").Value = 24508
$ws.Cells.Item(69, "C").Value = 0
$ws.Cells.Item(69, "D").Value = 7971
$ws.Cells.Item(69, "E").Value = 16282
$ws.Cells.Item(69, "F").Value = 0
$ws.Cells.Item(69, "G").Value = 0
$ws.Cells.Item(69, "H").Value = 255

$ws.Cells.Item(70, 1).Value = "Nepal"
$ws.Cells.Item(70, "B").Value = 24432
$ws.Cells.Item(70, "C").Value = 484
$ws.Cells.Item(70, "D").Value = 16728
$ws.Cells.Item(70, "E").Value = 7613
$ws.Cells.Item(70, "F").Value = 0
$ws.Cells.Item(70, "G").Value = 8
$ws.Cells.Item(70, "H").Value = 91

$ws.Cells.Item(125, 1).Value = "Mozambique"
$ws.Cells.Item(125, "B").Value = 2559
$ws.Cells.Item(125, "C").Value = 78
$ws.Cells.Item(125, "D").Value = 951
$ws.Cells.Item(125, "E").Value = 1589
$ws.Cells.Item(125, "F").Value = 0
$ws.Cells.Item(125, "G").Value = 2
$ws.Cells.Item(125, "H").Value = 19

$ws.Cells.Item(126, 1).Value = "Surinam"
$ws.Cells.Item(126, "B").Value = 2559
$ws.Cells.Item(126, "C").Value = 0
$ws.Cells.Item(126, "D").Value = 1712
$ws.Cells.Item(126, "E").Value = 808
$ws.Cells.Item(126, "F").Value = 0
$ws.Cells.Item(126, "G").Value = 0
$ws.Cells.Item(126, "H").Value = 39

$ws.Cells.Item(213, 1).Value = "Montserrat"
$ws.Cells.Item(213, "B").Value = 13
$ws.Cells.Item(213, "C").Value = 0
$ws.Cells.Item(213, "D").Value = 12
$ws.Cells.Item(213, "E").Value = 0
$ws.Cells.Item(213, "F").Value = 0
$ws.Cells.Item(213, "G").Value = 0
$ws.Cells.Item(213, "H").Value = 1

$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, "B").Value = 13
$ws.Cells.Item(214, "C").Value = 0
$ws.Cells.Item(214, "D").Value = 13
$ws.Cells.Item(214, "E").Value = 0
$ws.Cells.Item(214, "F").Value = 0
$ws.Cells.Item(214, "G").Value = 0
$ws.Cells.Item(214, "H").Value = 0

# --- Rows with refreshed statistics only (country unchanged) ---
$ws.Cells.Item(4, "B").Value = 5327918
$ws.Cells.Item(4, "C").Value = 21961
$ws.Cells.Item(4, "D").Value = 2778502
$ws.Cells.Item(4, "E").Value = 2381026
$ws.Cells.Item(4, "F").Value = 0
$ws.Cells.Item(4, "G").Value = 645
$ws.Cells.Item(4, "H").Value = 168390

$ws.Cells.Item(5, "B").Value = 3123109
$ws.Cells.Item(5, "C").Value = 10716
$ws.Cells.Item(5, "D").Value = 2243124
$ws.Cells.Item(5, "E").Value = 776564
$ws.Cells.Item(5, "F").Value = 0
$ws.Cells.Item(5, "G").Value = 322
$ws.Cells.Item(5, "H").Value = 103421

$ws.Cells.Item(6, "B").Value = 2395417
$ws.Cells.Item(6, "C").Value = 67012
$ws.Cells.Item(6, "D").Value = 1695860
$ws.Cells.Item(6, "E").Value = 652419
$ws.Cells.Item(6, "F").Value = 0
$ws.Cells.Item(6, "G").Value = 950
$ws.Cells.Item(6, "H").Value = 47138

$ws.Cells.Item(12, "B").Value = 378168
$ws.Cells.Item(12, "C").Value = 1552
$ws.Cells.Item(12, "D").Value = 351419
$ws.Cells.Item(12, "E").Value = 16544
$ws.Cells.Item(12, "F").Value = 0
$ws.Cells.Item(12, "G").Value = 27
$ws.Cells.Item(12, "H").Value = 10205

$ws.Cells.Item(23, "B").Value = 206696
$ws.Cells.Item(23, "C").Value = 2524
$ws.Cells.Item(23, "D").Value = 83237
$ws.Cells.Item(23, "E").Value = 93088
$ws.Cells.Item(23, "F").Value = 0
$ws.Cells.Item(23, "G").Value = 17
$ws.Cells.Item(23, "H").Value = 30371

$ws.Cells.Item(33, "B").Value = 87878
$ws.Cells.Item(33, "C").Value = 1285
$ws.Cells.Item(33, "D").Value = 61625
$ws.Cells.Item(33, "E").Value = 25614
$ws.Cells.Item(33, "F").Value = 0
$ws.Cells.Item(33, "G").Value = 17
$ws.Cells.Item(33, "H").Value = 639

$ws.Cells.Item(59, "B").Value = 36699
$ws.Cells.Item(59, "C").Value = 495
$ws.Cells.Item(59, "D").Value = 25627
$ws.Cells.Item(59, "E").Value = 9739
$ws.Cells.Item(59, "F").Value = 0
$ws.Cells.Item(59, "G").Value = 11
$ws.Cells.Item(59, "H").Value = 1333

$ws.Cells.Item(60, "B").Value = 36694
$ws.Cells.Item(60, "C").Value = 1499
$ws.Cells.Item(60, "D").Value = 25677
$ws.Cells.Item(60, "E").Value = 10461
$ws.Cells.Item(60, "F").Value = 0
$ws.Cells.Item(60, "G").Value = 23
$ws.Cells.Item(60, "H").Value = 556

$ws.Cells.Item(99, "B").Value = 6817
$ws.Cells.Item(99, "C").Value = 141
$ws.Cells.Item(99, "D").Value = 3552
$ws.Cells.Item(99, "E").Value = 3057
$ws.Cells.Item(99, "F").Value = 0
$ws.Cells.Item(99, "G").Value = 3
$ws.Cells.Item(99, "H").Value = 208

$ws.Cells.Item(105, "B").Value = 5223
$ws.Cells.Item(105, "C").Value = 0
$ws.Cells.Item(105, "D").Value = 2849
$ws.Cells.Item(105, "E").Value = 2353
$ws.Cells.Item(105, "F").Value = 0
$ws.Cells.Item(105, "G").Value = 1
$ws.Cells.Item(105, "H").Value = 21

$ws.Cells.Item(122, "B").Value = 2881
$ws.Cells.Item(122, "C").Value = 1
$ws.Cells.Item(122, "D").Value = 2638
$ws.Cells.Item(122, "E").Value = 232
$ws.Cells.Item(122, "F").Value = 0
$ws.Cells.Item(122, "G").Value = 0
$ws.Cells.Item(122, "H").Value = 11

$ws.Cells.Item(156, "B").Value = 963
$ws.Cells.Item(156, "C").Value = 0
$ws.Cells.Item(156, "D").Value = 855
$ws.Cells.Item(156, "E").Value = 56
$ws.Cells.Item(156, "F").Value = 0
$ws.Cells.Item(156, "G").Value = 0
$ws.Cells.Item(156, "H").Value = 52

$ws.Cells.Item(159, "B").Value = 880
$ws.Cells.Item(159, "C").Value = 14
$ws.Cells.Item(159, "D").Value = 403
$ws.Cells.Item(159, "E").Value = 460
$ws.Cells.Item(159, "F").Value = 0
$ws.Cells.Item(159, "G").Value = 1
$ws.Cells.Item(159, "H").Value = 17

